$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "WOMENS-JEWELRY-BUNDLEM" row (row 2) entirely; the rows below
# shift up so SONY-PS3-BUNDLEM moves to A2 and MICROSOFT-XBOX360-BUNDLEM
# moves to A3 (and the old last row disappears).
$ws.Rows.Item(2).Delete()

# The row that is now last (old A4/"MICROSOFT-XBOX360-BUNDLEM", now A3) used
# a wrapped-text style; keep that formatting intact on its new home.
$ws.Range("A3").WrapText = $true

# Move the active selection to A2, matching the post-edit selection.
$ws.Range("A2").Select()
